$wb = $excel.ActiveWorkbook

# Work on the "settings" sheet (3rd sheet in the workbook)
$settings = $wb.Worksheets.Item("settings")

# Add a "version" column with a value of 1, controlling the form version
$settings.Range("C1").Value = "version"
$settings.Range("C2").Value = 1

# Make the settings sheet the active sheet/tab (it receives focus for this edit)
$settings.Activate()
$settings.Range("C3").Select()
